# Automatische sync: 2025-06-17 13:57:41
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Add new log rows 11-14
$ws.Range("A11").Value = 'Re: Re: Wat zijn jullie openingstijden?'
$ws.Range("B11").Value = 'mailmind.test@zohomail.eu'
$ws.Range("C11").Value = 'Beste klant,
Hartelijk dank voor uw interesse. Onze openingstijden zijn maandag t/m vrijdag van 9:00 tot 18:00 uur en zaterdag van 10:00 tot 16:00 uur. Op zondag zijn wij gesloten. Voor verdere vragen kunt u altijd contact met ons opnemen.
Met vriendelijke groet,
[Naam van het bedrijf]'
$ws.Range("D11").Value = 'Informatieaanvraag'
$ws.Range("E11").Value = 'Geachte klant,
Dank voor uw interesse. Wij zijn geopend op maandag t/m vrijdag van 9:00-18:00 en zaterdag van 10:00-16:00. Op zondag zijn wij gesloten. Voor vragen zijn wij bereikbaar via info@bedrijfsnaam.nl of telefonisch op [telefoonnummer]. 
Met vriendelijke groet,
[Naam van het bedrijf]'
$ws.Range("F11").Value = '2025-06-17 12:59:24'
$ws.Range("G11").Value = 'Ja'

$ws.Range("A12").Value = 'Vragen over samenwerking'
$ws.Range("B12").Value = 'mailmind.test@zohomail.eu'
$ws.Range("C12").Value = 'Kunnen we samenwerken aan een nieuw project?'
$ws.Range("D12").Value = 'Overig'
$ws.Range("F12").Value = '2025-06-17 12:59:24'
$ws.Range("G12").Value = 'Nee'

$ws.Range("A13").Value = 'Sollicitatie marketingfunctie'
$ws.Range("B13").Value = 'mailmind.test@zohomail.eu'
$ws.Range("C13").Value = 'Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV.'
$ws.Range("D13").Value = 'Overig'
$ws.Range("F13").Value = '2025-06-17 12:59:25'
$ws.Range("G13").Value = 'Nee'

$ws.Range("A14").Value = 'Re: Re: Re: Wat zijn jullie openingstijden?'
$ws.Range("B14").Value = 'mailmind.test@zohomail.eu'
$ws.Range("C14").Value = 'Geachte klant,
Dank voor uw interesse. Wij zijn geopend op maandag t/m vrijdag van 9:00-18:00 en zaterdag van 10:00-16:00. Op zondag zijn wij gesloten. Voor vragen zijn wij bereikbaar via info@bedrijfsnaam.nl of telefonisch op [telefoonnummer]. 
Met vriendelijke groet,
[Naam van het bedrijf]'
$ws.Range("D14").Value = 'Informatieaanvraag'
$ws.Range("E14").Value = 'Geachte klant,
Bedankt voor uw interesse. Onze openingstijden zijn ma t/m vr van 9:00-18:00 en za van 10:00-16:00. Op zondag zijn we gesloten. Voor vragen zijn we bereikbaar via info@bedrijfsnaam.nl of telefonisch op [telefoonnummer].
Met vriendelijke groet,
[Naam van het bedrijf]'
$ws.Range("F14").Value = '2025-06-17 13:29:35'
$ws.Range("G14").Value = 'Ja'

# Expand conditional formatting ranges to include new rows
$ws.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D14"))
$ws.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G14"))

# Update Dashboard summary counts
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 6
$dash.Range("A3").Value = "Overig"
$dash.Range("B3").Value = 3
$dash.Range("A4").Value = "Afmelding"
$dash.Range("B4").Value = 2
